# Added more debug parameters
# Updates the "Translation" sheet of the TouchGFX texts workbook:
#  - Row 31 (STATE debug param) text is shortened from "STATE: <value>" to "S: <value>"
#  - Five new debug-parameter rows are appended (rows 33-37): B (bytes), M (memory),
#    R (fps reset) etc., mirroring the existing SingleUseIdNN / Medium / Left / LTR pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- Row 31: shorten existing "STATE: <value>" text to "S: <value>" ---
$ws.Range("F31").Value = "S: <value>"

# --- Row 33: SingleUseId51 / Medium / Left / LTR / "B: <value>" ---
$ws.Range("B33").Value = "SingleUseId51"
$ws.Range("C33").Value = "Medium"
$ws.Range("D33").Value = "Left"
$ws.Range("E33").Value = "LTR"
$ws.Range("F33").Value = "B: <value>"

# --- Row 34: SingleUseId52 / Medium / Left / LTR / "0" (stored as text) ---
$ws.Range("B34").Value = "SingleUseId52"
$ws.Range("C34").Value = "Medium"
$ws.Range("D34").Value = "Left"
$ws.Range("E34").Value = "LTR"
$ws.Range("F34").NumberFormat = "@"
$ws.Range("F34").Value = "0"
$ws.Range("F34").Style = "Normal"

# --- Row 35: SingleUseId53 / Medium / Left / LTR / "M: <value>" ---
$ws.Range("B35").Value = "SingleUseId53"
$ws.Range("C35").Value = "Medium"
$ws.Range("D35").Value = "Left"
$ws.Range("E35").Value = "LTR"
$ws.Range("F35").Value = "M: <value>"

# --- Row 36: SingleUseId54 / Medium / Left / LTR / "0" (stored as text) ---
$ws.Range("B36").Value = "SingleUseId54"
$ws.Range("C36").Value = "Medium"
$ws.Range("D36").Value = "Left"
$ws.Range("E36").Value = "LTR"
$ws.Range("F36").NumberFormat = "@"
$ws.Range("F36").Value = "0"
$ws.Range("F36").Style = "Normal"

# --- Row 37: SingleUseId55 / Medium / Center / LTR / "R" ---
$ws.Range("B37").Value = "SingleUseId55"
$ws.Range("C37").Value = "Medium"
$ws.Range("D37").Value = "Center"
$ws.Range("E37").Value = "LTR"
$ws.Range("F37").Value = "R"
